# Atualização da base de dados
# Updates "Inscritos" (E), "Pagos" (F) and "Inscrições homologadas" (H)
# columns for several rows of the "Inscricoes" worksheet/table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 56

$ws.Range("E4").Value = 20

$ws.Range("E10").Value = 172
$ws.Range("F10").Value = 67
$ws.Range("H10").Value = 67

$ws.Range("E11").Value = 117
$ws.Range("F11").Value = 54
$ws.Range("H11").Value = 54

$ws.Range("E12").Value = 177
$ws.Range("F12").Value = 85
$ws.Range("H12").Value = 85

$ws.Range("E13").Value = 58
$ws.Range("F13").Value = 24
$ws.Range("H13").Value = 24

$ws.Range("E14").Value = 55
$ws.Range("F14").Value = 23
$ws.Range("H14").Value = 23

$ws.Range("E15").Value = 74

$ws.Range("E17").Value = 33
$ws.Range("F17").Value = 14
$ws.Range("H17").Value = 14

$ws.Range("E18").Value = 27

$ws.Range("E20").Value = 44

$ws.Range("E21").Value = 59

$ws.Range("E22").Value = 73

$ws.Range("E23").Value = 76
$ws.Range("F23").Value = 30
$ws.Range("H23").Value = 30

$ws.Range("E24").Value = 84
$ws.Range("F24").Value = 36
$ws.Range("H24").Value = 36

$ws.Range("E26").Value = 43

$ws.Range("E27").Value = 112

$ws.Range("E28").Value = 74
$ws.Range("F28").Value = 19
$ws.Range("H28").Value = 19

$ws.Range("E29").Value = 73

$ws.Range("E30").Value = 75
$ws.Range("F30").Value = 35
$ws.Range("H30").Value = 35

$ws.Range("E31").Value = 34

$ws.Range("E32").Value = 78
$ws.Range("F32").Value = 37
$ws.Range("H32").Value = 37

$ws.Range("E33").Value = 107

$ws.Range("E34").Value = 82

$ws.Range("E35").Value = 52

$ws.Range("E36").Value = 23

$ws.Range("E37").Value = 54

$ws.Range("E38").Value = 39

$ws.Range("E39").Value = 85

$ws.Range("E41").Value = 145
$ws.Range("F41").Value = 49
$ws.Range("H41").Value = 49

$ws.Range("E42").Value = 126

$ws.Range("E44").Value = 116
$ws.Range("F44").Value = 51
$ws.Range("H44").Value = 51

$ws.Range("E45").Value = 49

$ws.Range("E46").Value = 90

$ws.Range("E47").Value = 164
$ws.Range("F47").Value = 63
$ws.Range("H47").Value = 63

$ws.Range("E48").Value = 80

$ws.Range("E49").Value = 89
$ws.Range("F49").Value = 39
$ws.Range("H49").Value = 39

$ws.Range("E50").Value = 66
$ws.Range("F50").Value = 26
$ws.Range("H50").Value = 26

$ws.Range("E51").Value = 72
$ws.Range("F51").Value = 29
$ws.Range("H51").Value = 29
